$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.550.76"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "2.281.97"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  +1.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.47%  "

$ws.Range("E9").Value = "  -1.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.974"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").Value = "2.629.55"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "2.308.87"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").Value = "42.502.81"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("E20").Value = "  -1.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("E28").Value = "  -1.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0858"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.71%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.24%  "

$ws.Range("E36").Value = "  -5.12%  "

$ws.Range("E37").Value = "  -2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0345"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "

$ws.Range("E40").Value = "  -3.72%  "

$ws.Range("E41").Value = "  +5.05%  "

$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.55%  "

$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.225"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").Value = "1.711.68"
$ws.Range("E47").Value = "  +7.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "79.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.17%  "

